# "Updated Excel test data"
#
# LoginTest!A2:G2 held a stale test fixture (Rahul / 8927342 /
# trainer@way2automation.com / India / Delhi / rahularora1985 / lsajdfksf).
# Replace it with the new fixture row, keep the county (India) as-is, widen
# column B to fit the longer phone number, and leave LoginTest as the
# selected/active sheet with D3 highlighted (NewCarsTest was the active tab
# before the edit).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginTest")

# New login fixture values (county/India stays the same).
$ws1.Range("A2").Value = "Mayank"
$ws1.Range("B2").Value = 3564684635
$ws1.Range("C2").Value = "mayank@mail.com"
$ws1.Range("E2").Value = "Lucknow"
$ws1.Range("F2").Value = "user"
$ws1.Range("G2").Value = "pass"

# Make LoginTest the active/selected sheet (it was NewCarsTest before) and
# move the selection to D3.
$ws1.Activate() | Out-Null
$ws1.Range("D3").Select() | Out-Null

# Column B (phoneno) is now wider to fit "3564684635".
$ws1.Columns.Item(2).ColumnWidth = 10.14
